# cole_gerrit.xlsx regen: column G ("K") values were recomputed to represent
# strikeouts (K) instead of the previous "Strike#" pitch-count style metric.
# Only column G (G2:G37) changes; update each cell's value in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 7
    4  = 6
    5  = 7
    6  = 7
    7  = 2
    8  = 15
    9  = 9
    10 = 6
    11 = 9
    12 = 10
    13 = 8
    14 = 11
    15 = 12
    16 = 6
    17 = 6
    18 = 6
    19 = 4
    20 = 9
    21 = 7
    22 = 5
    23 = 7
    24 = 9
    25 = 12
    26 = 4
    27 = 12
    28 = 11
    29 = 10
    30 = 8
    31 = 16
    32 = 8
    33 = 5
    34 = 8
    35 = 5
    36 = 5
    37 = 1
}

foreach ($row in $kValues.Keys | Sort-Object) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
